$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $escaped = $val.Replace('"', '""')
    $c = $ws.Range($addr)
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue $ws "D2" '30.586.19'
Set-TextValue $ws "E2" '  -0.01%  '
Set-TextValue $ws "D3" '1.919.11'
Set-TextValue $ws "E3" '  -0.17%  '
Set-TextValue $ws "E4" '  +0.15%  '
Set-TextValue $ws "D5" '246.22'
Set-TextValue $ws "E5" '  -0.29%  '
Set-TextValue $ws "D6" '1.005'
Set-TextValue $ws "E6" '  +0.46%  '
Set-TextValue $ws "D7" '0.4866'
Set-TextValue $ws "E7" '  +3.23%  '
Set-TextValue $ws "D8" '0.2892'
Set-TextValue $ws "E8" '  -0.20%  '
Set-TextValue $ws "D9" '0.06714'
Set-TextValue $ws "E9" '  -1.10%  '
Set-TextValue $ws "D10" '110.92'
Set-TextValue $ws "E10" '  +5.48%  '
Set-TextValue $ws "D11" '19.31'
Set-TextValue $ws "E11" '  +4.88%  '
Set-TextValue $ws "D12" '1.924.83'
Set-TextValue $ws "E12" '  +0.21%  '
Set-TextValue $ws "D13" '0.07612'
Set-TextValue $ws "E13" '  -1.15%  '
Set-TextValue $ws "D14" '5.332'
Set-TextValue $ws "E14" '  +0.75%  '
Set-TextValue $ws "D15" '0.6705'
Set-TextValue $ws "E15" '  -0.56%  '
Set-TextValue $ws "D16" '295.90'
Set-TextValue $ws "E16" '  +2.04%  '
Set-TextValue $ws "D17" '30.537.57'
Set-TextValue $ws "E17" '  -0.24%  '
Set-TextValue $ws "E18" '  +1.01%  '
Set-TextValue $ws "D19" '1.006'
Set-TextValue $ws "E19" '  +0.51%  '
Set-TextValue $ws "B20" 'ShibaInu'
Set-TextValue $ws "C20" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws "D20" '0.000007551'
Set-TextValue $ws "E20" '  -0.83%  '
Set-TextValue $ws "B21" 'Uniswap'
Set-TextValue $ws "C21" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws "D21" '5.529'
Set-TextValue $ws "E21" '  +1.38%  '
Set-TextValue $ws "D22" '2.172.65'
Set-TextValue $ws "E22" '  +0.40%  '
Set-TextValue $ws "D23" '1.001'
Set-TextValue $ws "E23" '  +0.04%  '
Set-TextValue $ws "D24" '6.498'
Set-TextValue $ws "E24" '  +2.52%  '
Set-TextValue $ws "D25" '9.456'
Set-TextValue $ws "E25" '  +0.51%  '
Set-TextValue $ws "D26" '164.62'
Set-TextValue $ws "E26" '  -2.04%  '
Set-TextValue $ws "E27" '  -3.44%  '
Set-TextValue $ws "D28" '2.095'
Set-TextValue $ws "E28" '  -1.28%  '
Set-TextValue $ws "E29" '  -0.55%  '
Set-TextValue $ws "D30" '1.453'
Set-TextValue $ws "E30" '  +6.23%  '
Set-TextValue $ws "D31" '4.145'
Set-TextValue $ws "E31" '  -0.92%  '
Set-TextValue $ws "D32" '4.036'
Set-TextValue $ws "E32" '  -2.64%  '
Set-TextValue $ws "D33" '0.05033'
Set-TextValue $ws "E33" '  -0.53%  '
Set-TextValue $ws "D34" '0.7392'
Set-TextValue $ws "E34" '  -0.40%  '
Set-TextValue $ws "D35" '1.141'
Set-TextValue $ws "E35" '  -1.17%  '
Set-TextValue $ws "D36" '0.9997'
Set-TextValue $ws "E36" '  +0.02%  '
Set-TextValue $ws "D37" '2.724'
Set-TextValue $ws "E37" '  -0.88%  '
Set-TextValue $ws "D38" '0.02028'
Set-TextValue $ws "E38" '  -2.86%  '
Set-TextValue $ws "D39" '2.687'
Set-TextValue $ws "E39" '  -0.25%  '
Set-TextValue $ws "D40" '110.77'
Set-TextValue $ws "E40" '  -0.35%  '
Set-TextValue $ws "D41" '2.018'
Set-TextValue $ws "E41" '  -2.07%  '
Set-TextValue $ws "D42" '0.4416'
Set-TextValue $ws "E42" '  +1.52%  '
Set-TextValue $ws "D43" '0.8666'
Set-TextValue $ws "E43" '  -1.75%  '
Set-TextValue $ws "B44" 'FraxShare'
Set-TextValue $ws "C44" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws "D44" '5.850'
Set-TextValue $ws "E44" '  -0.71%  '
Set-TextValue $ws "B45" 'Aave'
Set-TextValue $ws "C45" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws "D45" '70.82'
Set-TextValue $ws "E45" '  +5.33%  '
Set-TextValue $ws "D46" '1.007'
Set-TextValue $ws "E46" '  +0.69%  '
Set-TextValue $ws "D47" '7.237'
Set-TextValue $ws "E47" '  -0.08%  '
Set-TextValue $ws "D48" '48.34'
Set-TextValue $ws "E48" '  +0.89%  '
Set-TextValue $ws "D49" '9.145'
Set-TextValue $ws "E49" '  -1.71%  '
Set-TextValue $ws "D50" '0.1230'
Set-TextValue $ws "E50" '  -0.32%  '
Set-TextValue $ws "D51" '0.2534'
Set-TextValue $ws "E51" '  +4.29%  '

$ws.Application.CutCopyMode = $false

